# Developer Guide Section 3 update:
#
# The sequence diagram's class still shows the old "AddressBook" name in
# four places (":AddressBookParser", "undoAddressBook()",
# ":VersionedAddressBook" and "resetData(ReadOnlyAddressBook)"). Rename it
# to "ForumBook" by replacing the "Address" substring with "Forum" in each
# shape's text, leaving the rest of the label untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tf = $sh.TextFrame
        if ($tf.HasText) {
            $tr = $tf.TextRange
            $fullText = $tr.Text
            $pos = $fullText.IndexOf("Address")
            if ($pos -ge 0) {
                $sub = $tr.Characters($pos + 1, 7)
                $sub.Text = "Forum"
            }
        }
    }
}
